# JavaStudyRecord.pptx edit:
#  1. Collapse the multi-run hyperlink captions on slide 3 into single
#     runs per paragraph (the author had split each caption into many
#     <a:r> runs sharing identical formatting; they are merged back into
#     one run per paragraph, keeping the first run's formatting/hlink).
#  2. Remove the 4th slide ("Idea Tool Code Comments") from the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# --- Paragraph 3 (hyperlink rId2): 4 runs -> 1 run ------------------------
$para = $tr.Paragraphs(3, 1)
$para.Runs(1, 1).Text = "       Java集合源码系列之HashMap添加元素的流程"
for ($i = 0; $i -lt 3; $i++) {
    $para.Runs(2, 1).Text = ""
}

# --- Paragraph 4 (hyperlink rId3): leading run + 10 runs -> 1 run --------
$para = $tr.Paragraphs(4, 1)
$para.Runs(2, 1).Text = "Java基础篇：什么是hashCode 以及 hashCode()与equals()的联系"
for ($i = 0; $i -lt 9; $i++) {
    $para.Runs(3, 1).Text = ""
}

# --- Paragraph 5 (hyperlink rId4): leading run + 4 runs -> 1 run ---------
$para = $tr.Paragraphs(5, 1)
$para.Runs(2, 1).Text = "Java基础之hashcode剖析"
for ($i = 0; $i -lt 3; $i++) {
    $para.Runs(3, 1).Text = ""
}

# --- Paragraph 6 (hyperlink rId5): leading run + 4 runs -> 1 run ---------
$para = $tr.Paragraphs(6, 1)
$para.Runs(2, 1).Text = "深入理解 Java 中的 hashCode"
for ($i = 0; $i -lt 3; $i++) {
    $para.Runs(3, 1).Text = ""
}

# Re-running the text edits above shrinks the autofit shape (wrap="none",
# spAutoFit) because it now spans fewer rendered lines at the default
# measurement; restore the original autofit height (EMU 1991880 == this
# many points) so the shape geometry is unaffected by the text changes.
$shp.Height = 156.841

# --- Remove the 4th slide -------------------------------------------------
$p.Slides.Item(4).Delete()
